$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "varchar(50)"
$ws.Range("C1").Value = "varchar(50)"
$ws.Range("P1").Value = "varchar(8)"
$ws.Range("Q1").Value = "varchar(50)"
$ws.Range("R1").Value = "varchar(50)"

$ws.Range("A2").Value = "(id)입력x"
$ws.Range("B2").Value = "공장코드"
$ws.Range("C2").Value = "bom코드"
$ws.Range("D2").Value = "기초재고수량"
$ws.Range("E2").Value = "기초재고금액"
$ws.Range("F2").Value = "생산입고수량"
$ws.Range("G2").Value = "생산입고금액"
$ws.Range("H2").Value = "판매출고수량"
$ws.Range("I2").Value = "판매출고금액"
$ws.Range("J2").Value = "LOSS출고수량"
$ws.Range("K2").Value = "LOSS출고금액"
$ws.Range("L2").Value = "개발출고수량"
$ws.Range("M2").Value = "개발출고금액"
$ws.Range("N2").Value = "기말재고수량"
$ws.Range("O2").Value = "기말재고금액"
$ws.Range("P2").Value = "년월"
$ws.Range("Q2").Value = "버젼코드"
$ws.Range("R2").Value = "계정코드"
